# Generate Report for handoff
# b.md is now "Ready for handoff": update status cells, the zh-cn/de-de
# handoff-file hyperlink display text, and the handoff datetimes.

$wb = $excel.ActiveWorkbook

$missing = [System.Reflection.Missing]::Value

# ---------------------------------------------------------------------------
# Sheet "Overview": b.md row (row 3) status -> "Ready for handoff"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------------
# Language sheets (zh-cn / de-de): the host's Hyperlinks collection only
# supports Add (always appends) and a bulk Delete (wipes every hyperlink on
# the sheet) -- there is no in-place update. Reading back .Address /
# .TextToDisplay on existing links also comes back empty, so the known
# addresses/labels are supplied explicitly here (they match the workbook's
# existing relationships) and the collection is rebuilt from scratch with
# the one cell (C3) that changes getting its new display text.
# ---------------------------------------------------------------------------

function Rebuild-Hyperlinks {
    param($ws, $links)

    $ws.Hyperlinks.Delete()
    foreach ($link in $links) {
        $ws.Hyperlinks.Add($ws.Range($link.ref), $link.addr, $missing, $missing, $link.display)
    }
}

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "2016-01-22 02:25:33"

$zhLinks = @(
    @{ ref = "A2"; display = "a.md.md"; addr = "https://github.com/OpenLocalizationTest/oltest/blob/be2c191677cfa29138830218249bff4ab8a8ce01/e2e/a.md.md" },
    @{ ref = "C2"; display = "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"; addr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2eeab7f5ea2a7a51b9833eb39c75cc233ab2e2eb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf" },
    @{ ref = "E2"; display = "a.md.md"; addr = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/dacf8d24b45a5307190e9b70016c798f90539c6a/e2e/a.md.md" },
    @{ ref = "F2"; display = "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"; addr = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/63f69aaf08c8f7cb4902fcd9f84d3ee950524fa0/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf" },
    @{ ref = "A3"; display = "b.md.md"; addr = "https://github.com/OpenLocalizationTest/oltest/blob/be2c191677cfa29138830218249bff4ab8a8ce01/e2e/b.md.md" },
    @{ ref = "C3"; display = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"; addr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2eeab7f5ea2a7a51b9833eb39c75cc233ab2e2eb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf" },
    @{ ref = "E3"; display = "a.md.md"; addr = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/dacf8d24b45a5307190e9b70016c798f90539c6a/e2e/a.md.md" },
    @{ ref = "F3"; display = "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"; addr = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/63f69aaf08c8f7cb4902fcd9f84d3ee950524fa0/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf" },
    @{ ref = "A4"; display = ".localization-config"; addr = "https://github.com/OpenLocalizationTest/oltest/blob/be2c191677cfa29138830218249bff4ab8a8ce01/.localization-config" }
)
Rebuild-Hyperlinks $wsZh $zhLinks

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "2016-01-22 02:25:47"

$deLinks = @(
    @{ ref = "A2"; display = "a.md.md"; addr = "https://github.com/OpenLocalizationTest/oltest/blob/be2c191677cfa29138830218249bff4ab8a8ce01/e2e/a.md.md" },
    @{ ref = "C2"; display = "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"; addr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8dfa17969aa07f13d578e9ac61f7dd1455519fe5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf" },
    @{ ref = "E2"; display = "a.md.md"; addr = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4b90a1d8f7348ce08f802f07fb05c6f5dc367628/e2e/a.md.md" },
    @{ ref = "F2"; display = "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"; addr = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9fc777158e94625b39270cb9267403a43c98b5c3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf" },
    @{ ref = "A3"; display = "b.md.md"; addr = "https://github.com/OpenLocalizationTest/oltest/blob/be2c191677cfa29138830218249bff4ab8a8ce01/e2e/b.md.md" },
    @{ ref = "C3"; display = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"; addr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8dfa17969aa07f13d578e9ac61f7dd1455519fe5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf" },
    @{ ref = "E3"; display = "a.md.md"; addr = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4b90a1d8f7348ce08f802f07fb05c6f5dc367628/e2e/a.md.md" },
    @{ ref = "F3"; display = "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"; addr = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9fc777158e94625b39270cb9267403a43c98b5c3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf" },
    @{ ref = "A4"; display = ".localization-config"; addr = "https://github.com/OpenLocalizationTest/oltest/blob/be2c191677cfa29138830218249bff4ab8a8ce01/.localization-config" }
)
Rebuild-Hyperlinks $wsDe $deLinks
